# Update the "as_of_utc" timestamp (column AA, rows 2-26) on both the
# "Главные" and "Линейные" sheets from 2025-11-24 03:04:16 to
# 2025-11-24 07:04:22.

$wb = $excel.ActiveWorkbook
$newTimestamp = "2025-11-24 07:04:22"

for ($sheetIndex = 2; $sheetIndex -le 3; $sheetIndex++) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
